# "fixing vstu -> vsu"
#
# The link text "https://github.com/vitali-kviatkouski/vstu" on the
# "Links" slide is corrected to "https://github.com/vitali-kviatkouski/vsu".
# In the canonical OOXML this shows up as the single run being split into
# two runs (both keeping the same formatting / hyperlink):
#   <a:r>...<a:t>https://</a:t></a:r>
#   <a:r>...<a:t>github.com/vitali-kviatkouski/vsu</a:t></a:r>
# so we reproduce that here instead of just overwriting the whole string,
# which would otherwise collapse back into a single run.

$p = $ppt.ActivePresentation

$oldUrlText = "https://github.com/vitali-kviatkouski/vstu"
$newTailText = "github.com/vitali-kviatkouski/vsu"
$splitText = "https://"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text

        if ($fullText.IndexOf($oldUrlText) -lt 0) {
            continue
        }

        $startIdx = $fullText.IndexOf($oldUrlText)

        # Narrow down to just the run's text and trim it to the "https://" prefix.
        $runRange = $tr.Characters($startIdx + 1, $oldUrlText.Length)
        $runRange.Text = $splitText

        # Insert the corrected remainder right after it as its own run.
        $prefixRange = $tr.Characters($startIdx + 1, $splitText.Length)
        $tailRange = $prefixRange.InsertAfter($newTailText)

        # Re-asserting the (unchanged) font size forces the engine to keep
        # the insertion as a distinct run instead of merging it back into
        # the previous one, matching the two-run shape of the target XML.
        $tailRange.Font.Size = $runRange.Font.Size
    }
}
